# "Generate Report for Archive" - refresh the localization status report:
#  - flip the in-flight status text from "Ready for handoff" to "In Translation"
#  - narrow the Status column now that the shorter text fits

$wb = $excel.ActiveWorkbook

# ---- Update status text on every sheet that shows it ----------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# ---- Narrow the Status columns to match the new content -------------------
$overview.Columns("E").ColumnWidth = 12.5
$overview.Columns("F").ColumnWidth = 12.5
$zhcn.Columns("C").ColumnWidth = 12.5
$dede.Columns("C").ColumnWidth = 12.5
